$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Arts for Lawrence
$ws.Range("A3").Value = 33

# Row 7: Cumberland Arts, Inc.
$ws.Range("A7").Value = 22
$ws.Range("C7").Value = 'Cumberland Arts, Inc.'
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0

# Row 8: Downtown Indy, Inc.
$ws.Range("A8").Value = 23
$ws.Range("C8").Value = 'Downtown Indy, Inc.'
$ws.Range("D8").Value = 4.6
$ws.Range("E8").Value = 246

# Row 9: Festival Flea Market
$ws.Range("A9").Value = 21
$ws.Range("C9").Value = 'Festival Flea Market'
$ws.Range("D9").Value = 3.9
$ws.Range("E9").Value = 191

# Row 10: Free Stage
$ws.Range("A10").Value = 25
$ws.Range("C10").Value = 'Free Stage'
$ws.Range("D10").Value = 4.6
$ws.Range("E10").Value = 361

# Row 11: HSI Show Productions
$ws.Range("A11").Value = 14
$ws.Range("C11").Value = 'HSI Show Productions'
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0

# Row 12: Heartland Film
$ws.Range("A12").Value = 19
$ws.Range("C12").Value = 'Heartland Film'
$ws.Range("D12").Value = 4.3
$ws.Range("E12").Value = 14

# Row 13: Holy Rosary Church
$ws.Range("A13").Value = 24
$ws.Range("C13").Value = 'Holy Rosary Church'
$ws.Range("D13").Value = 4.8
$ws.Range("E13").Value = 210

# Row 14: Indiana Arts Building
$ws.Range("A14").Value = 18
$ws.Range("C14").Value = 'Indiana Arts Building'
$ws.Range("E14").Value = 3

# Row 15: Indianapolis Auto Show
$ws.Range("A15").Value = 13
$ws.Range("C15").Value = 'Indianapolis Auto Show'
$ws.Range("D15").Value = 3.1
$ws.Range("E15").Value = 51

# Row 16: Indianapolis Grapevine
$ws.Range("A16").Value = 27
$ws.Range("C16").Value = 'Indianapolis Grapevine'
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 0

# Row 17: Indianapolis Jewish Film Festival
$ws.Range("A17").Value = 16
$ws.Range("C17").Value = 'Indianapolis Jewish Film Festival'
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 0

# Row 18: Indianapolis Zoo
$ws.Range("A18").Value = 28
$ws.Range("C18").Value = 'Indianapolis Zoo'
$ws.Range("D18").Value = 4.5
$ws.Range("E18").Value = 14982

# Row 19: Indy Taco Fest
$ws.Range("A19").Value = 0
$ws.Range("C19").Value = 'Indy Taco Fest'
$ws.Range("D19").Value = 4.2
$ws.Range("E19").Value = 30

# Row 20: IndyFringe Theatre
$ws.Range("A20").Value = 32
$ws.Range("C20").Value = 'IndyFringe Theatre'
$ws.Range("D20").Value = 4.6
$ws.Range("E20").Value = 170

# Row 21: Irvington Halloween Festival
$ws.Range("A21").Value = 1
$ws.Range("C21").Value = 'Irvington Halloween Festival'
$ws.Range("D21").Value = 4.7
$ws.Range("E21").Value = 13

# Row 22: Lights Over Seminary
$ws.Range("A22").Value = 34
$ws.Range("C22").Value = 'Lights Over Seminary'
$ws.Range("D22").Value = 5
$ws.Range("E22").Value = 2

# Row 23: Lourdes Lyons Chili Cook-Off
$ws.Range("A23").Value = 12
$ws.Range("C23").Value = 'Lourdes Lyons Chili Cook-Off'
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0

# Row 24: Marion County Fairgrounds
$ws.Range("A24").Value = 29
$ws.Range("C24").Value = 'Marion County Fairgrounds'
$ws.Range("D24").Value = 4.2
$ws.Range("E24").Value = 1104

# Row 25: MasterWorks Festival
$ws.Range("A25").Value = 17
$ws.Range("C25").Value = 'MasterWorks Festival'
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0

# Row 26: Military Park
$ws.Range("A26").Value = 11
$ws.Range("C26").Value = 'Military Park'
$ws.Range("D26").Value = 4.6
$ws.Range("E26").Value = 1512

# Row 27: New Palestine Lions Club
$ws.Range("A27").Value = 31
$ws.Range("C27").Value = 'New Palestine Lions Club'
$ws.Range("D27").Value = 4.6
$ws.Range("E27").Value = 54

# Row 28: Nickel Plate District Amphitheater
$ws.Range("A28").Value = 30
$ws.Range("C28").Value = 'Nickel Plate District Amphitheater'
$ws.Range("D28").Value = 4.6
$ws.Range("E28").Value = 430

# Row 29: Spirit & Place Festival
$ws.Range("A29").Value = 10
$ws.Range("C29").Value = 'Spirit & Place Festival'

# Row 30: St Louis de Montfort Fall Festival
$ws.Range("A30").Value = 3
$ws.Range("C30").Value = 'St Louis de Montfort Fall Festival'
$ws.Range("D30").Value = 4
$ws.Range("E30").Value = 4

# Row 31: St. Christopher Mid-Summer Festival
$ws.Range("A31").Value = 4
$ws.Range("C31").Value = 'St. Christopher Mid-Summer Festival'
$ws.Range("D31").Value = 4.5
$ws.Range("E31").Value = 13

# Row 32: Stonycreek Farm Presents: HauntFest
$ws.Range("A32").Value = 26
$ws.Range("C32").Value = 'Stonycreek Farm Presents: HauntFest'
$ws.Range("D32").Value = 4
$ws.Range("E32").Value = 71

# Row 33: Talbot Street Art Fair
$ws.Range("A33").Value = 5
$ws.Range("C33").Value = 'Talbot Street Art Fair'
$ws.Range("E33").Value = 99

# Row 34: Viking Fest
$ws.Range("A34").Value = 6
$ws.Range("C34").Value = 'Viking Fest'
$ws.Range("D34").Value = 4.7
$ws.Range("E34").Value = 3

# Row 35: Waterman's Family Farm
$ws.Range("A35").Value = 15
$ws.Range("C35").Value = 'Waterman''s Family Farm'
$ws.Range("D35").Value = 4.4
$ws.Range("E35").Value = 669

# Row 36: indianapolis greek fest
$ws.Range("A36").Value = 2
$ws.Range("C36").Value = 'indianapolis greek fest'
$ws.Range("D36").Value = 4.3
$ws.Range("E36").Value = 19

# Remove now-obsolete trailing rows (41 -> 36 rows)
$ws.Range("A37:E41").EntireRow.Delete()
